# Auto-generated edit script: updates per diff mapping (Sheets/Gilgamesh_Profits.xlsx)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1364.95
$ws.Range("I15").Value = 1364.95
$ws.Range("K15").Value = 4094.85
$ws.Range("M15").Value = -3925.85

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()  # was -999028

$ws.Range("H55").Value = 586.1667
$ws.Range("J55").Value = 577
$ws.Range("L55").Value = 577
$ws.Range("N55").Value = -1005

$ws.Range("H107").Value = 589.4375
$ws.Range("I107").Value = 590.8570999999999
$ws.Range("J107").Value = 579.5
$ws.Range("K107").Value = 590.8570999999999
$ws.Range("L107").Value = 579.5
$ws.Range("M107").Value = 1329.1429
$ws.Range("N107").Value = -4419.5

$ws.Range("H112").Value = 2416.963
$ws.Range("J112").Value = 2440.6924
$ws.Range("L112").Value = 7322.0772
$ws.Range("N112").Value = -9538.0772

$ws.Range("H113").Value = 66669868
$ws.Range("I113").Value = 166668670
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 166668670
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = -166665416
$ws.Range("N113").Value = -10507

$ws.Range("H138").Value = 469002.97
$ws.Range("J138").Value = 682873.4399999999
$ws.Range("L138").Value = 2048620.32
$ws.Range("N138").Value = -2058900.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2118.4583
$ws.Range("J2").Value = 2099.8333
$ws.Range("L2").Value = 2099.8333
$ws.Range("N2").Value = -2325.8333

$ws.Range("H32").Value = 15599.474
$ws.Range("I32").Value = 9654.933999999999
$ws.Range("K32").Value = 9654.933999999999
$ws.Range("M32").Value = -9367.933999999999

$ws.Range("H37").Value = 59998
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()  # was -70584

$ws.Range("H74").Value = 138182.88
$ws.Range("I74").Value = 216206.73
$ws.Range("J74").Value = 2941.5334
$ws.Range("K74").Value = 216206.73
$ws.Range("L74").Value = 2941.5334
$ws.Range("M74").Value = -215332.73
$ws.Range("N74").Value = -4689.5334

$ws.Range("H77").Value = 138182.88
$ws.Range("I77").Value = 216206.73
$ws.Range("J77").Value = 2941.5334
$ws.Range("K77").Value = 1081033.65
$ws.Range("L77").Value = 14707.667
$ws.Range("M77").Value = -1076665.65
$ws.Range("N77").Value = -23443.667

$ws.Range("H110").Value = 3607.5715
$ws.Range("I110").Value = 2542.5833
$ws.Range("K110").Value = 2542.5833
$ws.Range("M110").Value = -497.5832999999998

$ws.Range("H116").Value = 2118.4583
$ws.Range("J116").Value = 2099.8333
$ws.Range("L116").Value = 2099.8333
$ws.Range("N116").Value = -6687.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2118.4583
$ws.Range("J3").Value = 2099.8333
$ws.Range("L3").Value = 2099.8333
$ws.Range("N3").Value = -2327.8333

$ws.Range("H99").Value = 69049.3
$ws.Range("I99").Value = 85324.30499999999
$ws.Range("K99").Value = 85324.30499999999
$ws.Range("M99").Value = -83826.30499999999

$ws.Range("H105").Value = 12383578
$ws.Range("I105").Value = 627215.2
$ws.Range("K105").Value = 627215.2
$ws.Range("M105").Value = -625468.2

$ws.Range("H107").Value = 1531.6923
$ws.Range("I107").Value = 1469.5555
$ws.Range("J107").Value = 1671.5
$ws.Range("K107").Value = 1469.5555
$ws.Range("L107").Value = 1671.5
$ws.Range("M107").Value = 450.4445000000001
$ws.Range("N107").Value = -5511.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14242.857
$ws.Range("I31").Value = 12500
$ws.Range("J31").Value = 14940
$ws.Range("K31").Value = 12500
$ws.Range("L31").Value = 14940
$ws.Range("M31").Value = -12205
$ws.Range("N31").Value = -15530

$ws.Range("H34").Value = 14242.857
$ws.Range("I34").Value = 12500
$ws.Range("J34").Value = 14940
$ws.Range("K34").Value = 12500
$ws.Range("L34").Value = 14940
$ws.Range("M34").Value = -12298
$ws.Range("N34").Value = -15344

$ws.Range("H58").Value = 3529.7742
$ws.Range("I58").Value = 2921.5
$ws.Range("J58").Value = 4635.727
$ws.Range("K58").Value = 2921.5
$ws.Range("L58").Value = 4635.727
$ws.Range("M58").Value = -2718.5
$ws.Range("N58").Value = -5041.727

$ws.Range("H63").Value = 80000
$ws.Range("J63").Value = 80000
$ws.Range("L63").Value = 80000
$ws.Range("N63").Value = -81372

$ws.Range("H66").Value = 80000
$ws.Range("J66").Value = 80000
$ws.Range("L66").Value = 240000
$ws.Range("N66").Value = -246864

$ws.Range("H132").Value = 4513.469
$ws.Range("I132").Value = 3920.5386
$ws.Range("J132").Value = 6825.9
$ws.Range("K132").Value = 11761.6158
$ws.Range("L132").Value = 20477.7
$ws.Range("M132").Value = -9231.6158
$ws.Range("N132").Value = -25537.7

$ws.Range("H136").Value = 3529.7742
$ws.Range("I136").Value = 2921.5
$ws.Range("J136").Value = 4635.727
$ws.Range("K136").Value = 8764.5
$ws.Range("L136").Value = 13907.181
$ws.Range("M136").Value = -6214.5
$ws.Range("N136").Value = -19007.181

$ws.Range("H141").Value = 884727.6
$ws.Range("J141").Value = 884727.6
$ws.Range("L141").Value = 884727.6
$ws.Range("N141").Value = -895087.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 52565628
$ws.Range("I4").Value = 66191988
$ws.Range("J4").Value = 6235999.5
$ws.Range("K4").Value = 198575964
$ws.Range("L4").Value = 18707998.5
$ws.Range("M4").Value = -198575852
$ws.Range("N4").Value = -18708222.5

$ws.Range("H98").Value = 397.875
$ws.Range("I98").Value = 442.6
$ws.Range("K98").Value = 1327.8
$ws.Range("M98").Value = 170.1999999999998

$ws.Range("H107").Value = 1346.5385
$ws.Range("J107").Value = 1552.5
$ws.Range("L107").Value = 4657.5
$ws.Range("N107").Value = -8497.5

$ws.Range("H113").Value = 6699.55
$ws.Range("I113").Value = 1997
$ws.Range("J113").Value = 6947.0527
$ws.Range("K113").Value = 5991
$ws.Range("L113").Value = 20841.1581
$ws.Range("M113").Value = -3821
$ws.Range("N113").Value = -25181.1581

$ws.Range("H121").Value = 200499.9
$ws.Range("I121").Value = 671
$ws.Range("J121").Value = 666767.3
$ws.Range("K121").Value = 2013
$ws.Range("L121").Value = 2000301.9
$ws.Range("M121").Value = -703
$ws.Range("N121").Value = -2002921.9

$ws.Range("H122").Value = 2305.75
$ws.Range("J122").Value = 2398.25
$ws.Range("L122").Value = 21584.25
$ws.Range("N122").Value = -26484.25

$ws.Range("H123").Value = 2166.6667
$ws.Range("I123").Value = 1678.5714
$ws.Range("K123").Value = 5035.7142
$ws.Range("M123").Value = -2585.7142

$ws.Range("H132").Value = 1841.2307
$ws.Range("I132").Value = 1652.9412
$ws.Range("K132").Value = 14876.4708
$ws.Range("M132").Value = -12346.4708

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 47749.25
$ws.Range("J62").Value = 47999
$ws.Range("L62").Value = 47999
$ws.Range("N62").Value = -49371

$ws.Range("H65").Value = 47749.25
$ws.Range("J65").Value = 47999
$ws.Range("L65").Value = 143997
$ws.Range("N65").Value = -150861

$ws.Range("H102").Value = 8942.964
$ws.Range("I102").Value = 1685.5714
$ws.Range("K102").Value = 1685.5714
$ws.Range("M102").Value = -63.57140000000004

$ws.Range("H132").Value = 3840.5
$ws.Range("I132").Value = 2898.6667
$ws.Range("K132").Value = 8696.000100000001
$ws.Range("M132").Value = -6166.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 30499.5
$ws.Range("J62").Value = 30499.5
$ws.Range("L62").Value = 30499.5
$ws.Range("N62").Value = -31747.5

$ws.Range("H63").Value = 49999
$ws.Range("J63").Value = 49999
$ws.Range("L63").Value = 49999
$ws.Range("N63").Value = -51497

$ws.Range("H65").Value = 30499.5
$ws.Range("J65").Value = 30499.5
$ws.Range("L65").Value = 91498.5
$ws.Range("N65").Value = -97738.5

$ws.Range("H66").Value = 49999
$ws.Range("J66").Value = 49999
$ws.Range("L66").Value = 149997
$ws.Range("N66").Value = -157485

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 16673380
$ws.Range("I122").Value = 5058.1665
$ws.Range("J122").Value = 83346670
$ws.Range("K122").Value = 15174.4995
$ws.Range("L122").Value = 250040010
$ws.Range("M122").Value = -12724.4995
$ws.Range("N122").Value = -250044910

$ws.Range("H123").Value = 69161.25
$ws.Range("J123").Value = 69161.25
$ws.Range("L123").Value = 69161.25
$ws.Range("N123").Value = -78961.25
